$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 76

# Column A holds a date-looking string ("2025-10-13") that must stay as
# literal text (matching every other row in this sheet), not be promoted
# to a real date serial by Excel's automatic type inference. Briefly mark
# the cell as Text before assigning, then restore the default "Normal"
# style so the cell ends up with no explicit formatting (same as its
# neighbours).
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025-10-13"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = "21:21:35"
$ws.Cells.Item($row, 3).Value = "1.00 EUR = 1,753.3027"
